$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 800 (shifts rows 800:841 down to 801:842,
# growing the used range from A1:D841 to A1:D842).
$ws.Rows.Item(800).Insert()

# Populate the newly inserted row with the new data point
# (2026/02/12, Thursday, hour 13, ranking 201).
# The leading apostrophe forces the date-looking text to be stored
# as literal text (matching the rest of column A) instead of being
# auto-converted into a real Excel date serial number; resetting the
# style back to "Normal" afterwards drops the quote-prefix flag that
# the apostrophe trick leaves behind, so the cell ends up as a plain
# string cell just like its neighbours.
$ws.Range("A800").Value = "'2026/02/12"
$ws.Range("A800").Style = "Normal"

$ws.Range("B800").Value = "木"

$ws.Range("C800").Value = 13
$ws.Range("D800").Value = 201
